$wb = $excel.ActiveWorkbook

# The first (existing) sheet - used as a style/format donor for the new sheet
$ws1 = $wb.Worksheets.Item(1)

# Add the new "StockTestData" sheet right after the existing sheet
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "StockTestData"

# --- Data ---
$newSheet.Range("A1").Value = "AddStockTest"

$newSheet.Range("A2").Value = "PortfolioName"
$newSheet.Range("B2").Value = "StockName"
$newSheet.Range("C2").Value = "StockQty"
$newSheet.Range("D2").Value = "StockPrice"
$newSheet.Range("E2").Value = "Comments"

$newSheet.Range("A3").Value = "My 2023 Stock"
$newSheet.Range("B3").Value = "ITC"
$newSheet.Range("C3").Value = 1000
$newSheet.Range("D3").Value = 420

# --- Formatting (reuse the existing section-title / header-row formats) ---
$ws1.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$ws1.Range("A2:D2").Copy()
$newSheet.Range("A2:D2").PasteSpecial(-4122)

$ws1.Range("D2").Copy()
$newSheet.Range("E2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Selection / active cell on the new sheet ---
$newSheet.Range("E11").Select()
